# Fills in "Programador # 2" (Rodrigo Daniel Salvatierra Morales - 2017288)
# and underlines the word "Programador" in the "Programador # 5" heading.
#
# Commit: "Se agrego Programador #2. (rsalvatierra-2017288)"

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two placeholder paragraphs that immediately follow the
# "Programador # 2" heading paragraph (name placeholder, then activity
# placeholder) by walking the Paragraphs collection instead of hard
# coding indices, so the script is resilient to minor structural drift.
# ------------------------------------------------------------------
$progTwoIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Programador # 2*") {
        $progTwoIndex = $i
        break
    }
}

if ($progTwoIndex -eq -1) {
    throw "Could not find the 'Programador # 2' paragraph."
}

$namePara = $d.Paragraphs.Item($progTwoIndex + 1)
$activityPara = $d.Paragraphs.Item($progTwoIndex + 2)

# --- Name placeholder: "(Escribir nombre completo y número de carnet)" ---
# Becomes two runs: "Rodrigo Daniel Salvatierra" + " Morales - 2017288"
$nameRange = $namePara.Range
[void]$nameRange.MoveEnd(1, -1) # exclude the paragraph mark
if ($nameRange.Text -eq "(Escribir nombre completo y número de carnet)") {
    $nameRange.Text = "Rodrigo Daniel Salvatierra"

    $tailRange = $namePara.Range
    [void]$tailRange.MoveEnd(1, -1)
    [void]$tailRange.Collapse(0)
    [void]$tailRange.InsertAfter(" Morales - 2017288")
    # Nudge formatting off/on so the inserted text is written out as its
    # own run instead of being silently merged back into the previous one.
    $tailRange.Bold = 1
    $tailRange.Bold = 0
}

# --- Activity placeholder: "(Escribir actividad asignada)" ---
$activityRange = $activityPara.Range
[void]$activityRange.MoveEnd(1, -1)
if ($activityRange.Text -eq "(Escribir actividad asignada)") {
    $activityRange.Text = "Creación de Base de datos con insercion de datos de cada entidad"
}

# ------------------------------------------------------------------
# "Programador # 5" -> underline just the word "Programador", leaving
# " # 5" without underline, split across two runs.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    [void]$full.MoveEnd(1, -1)
    if ($full.Text -eq "Programador # 5") {
        $wordRange = $full.Duplicate
        [void]$wordRange.SetRange($full.Start, $full.Start + 11) # "Programador" = 11 chars
        $wordRange.Font.Underline = 1
        break
    }
}
